# Resubmission based on mid-semester feedback
# Remove the empty/blank slide (sldId 262 -> Slides.Item(7)) that sits
# between the "Case 3/4" comparison slides and the final summary slide.
# Title 1 / Content Placeholder 2 on that slide have no body text, so it
# was cut from the deck; the slide formerly known as sldId 263 becomes
# the new slide 7.

$p = $ppt.ActivePresentation
$p.Slides.Item(7).Delete()
